# Add 5 new date columns (LM:LQ) to the mobility sheet, mirroring the
# existing last column (LL) for formatting, and fill in the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column LL's formatting (incl. the date number format used by row 1)
# into the new LM:LQ columns so the new cells pick up the same style index
# as the existing trailing date column instead of creating a new style.
$ws.Range("LL1:LL5").Copy($ws.Range("LM1:LQ5"))

# Row 1 holds date serials (continuing the daily sequence from LL1=44160).
$ws.Range("LM1").Value = 44161
$ws.Range("LN1").Value = 44162
$ws.Range("LO1").Value = 44163
$ws.Range("LP1").Value = 44164
$ws.Range("LQ1").Value = 44165

# Row 2 data values.
$ws.Range("LM2").Value = 37.71
$ws.Range("LN2").Value = 54.6
$ws.Range("LO2").Value = 51.5
$ws.Range("LP2").Value = 46.27
$ws.Range("LQ2").Value = 57.71

# Row 3 data values.
$ws.Range("LM3").Value = 29.8
$ws.Range("LN3").Value = 43.9
$ws.Range("LO3").Value = 46.86
$ws.Range("LP3").Value = 41
$ws.Range("LQ3").Value = 39.72

# Row 4 data values.
$ws.Range("LM4").Value = 37.08
$ws.Range("LN4").Value = 57.45
$ws.Range("LO4").Value = 56.84
$ws.Range("LP4").Value = 52.89
$ws.Range("LQ4").Value = 63.26

# Row 5 data values.
$ws.Range("LM5").Value = 22.23
$ws.Range("LN5").Value = 30
$ws.Range("LO5").Value = 28.43
$ws.Range("LP5").Value = 25.4
$ws.Range("LQ5").Value = 32.32
